$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = 'Volume 30   Number  33'
$ws.Range("C9").Value = 'Report Covering the Week  8/14/2023  Through  8/20/2023'

# Row 15
$ws.Range("D15").Value = '0'
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = '***.*'
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 83.333333333333
$ws.Range("N15").Value = -26.666666666666

# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 173
$ws.Range("J16").Value = 154
$ws.Range("K16").Value = 12.337662337662
$ws.Range("L16").Value = 33.076923076923
$ws.Range("M16").Value = -13.5
$ws.Range("N16").Value = -84.525939177102

# Row 17
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 116.666666666667
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 44
$ws.Range("H17").Value = 15.909090909090
$ws.Range("I17").Value = 308
$ws.Range("J17").Value = 335
$ws.Range("K17").Value = -8.059701492537
$ws.Range("L17").Value = 19.379844961240
$ws.Range("M17").Value = 66.486486486486
$ws.Range("N17").Value = -38.645418326693

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 47.058823529411
$ws.Range("I18").Value = 147
$ws.Range("J18").Value = 178
$ws.Range("K18").Value = -17.415730337078
$ws.Range("L18").Value = -9.815950920245
$ws.Range("M18").Value = -51.803278688524
$ws.Range("N18").Value = -89.309090909090

# Row 19
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 76
$ws.Range("G19").Value = 80
$ws.Range("H19").Value = -5
$ws.Range("I19").Value = 577
$ws.Range("J19").Value = 559
$ws.Range("K19").Value = 3.220035778175
$ws.Range("L19").Value = 35.764705882352
$ws.Range("M19").Value = 34.498834498834
$ws.Range("N19").Value = -16.133720930232

# Row 20
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 43
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 43.333333333333
$ws.Range("I20").Value = 253
$ws.Range("J20").Value = 198
$ws.Range("K20").Value = 27.777777777777
$ws.Range("L20").Value = 90.225563909774
$ws.Range("M20").Value = 42.937853107344
$ws.Range("N20").Value = -85.117647058823

# Row 21
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = 37.5
$ws.Range("F21").Value = 216
$ws.Range("G21").Value = 194
$ws.Range("H21").Value = 11.340206185567
$ws.Range("I21").Value = 1484
$ws.Range("J21").Value = 1445
$ws.Range("K21").Value = 2.698961937716
$ws.Range("L21").Value = 30.404217926186
$ws.Range("M21").Value = 12.851711026616
$ws.Range("N21").Value = -72.665315896113

# Row 22
$ws.Range("D22").Value = '0'
$ws.Range("N22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = '***.*'
$ws.Range("N22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 25
$ws.Range("K22").Value = 78.571428571428
$ws.Range("L22").Value = 31.578947368421
$ws.Range("M22").Value = 25

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 33.333333333333
$ws.Range("F23").Value = 26
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = 13.043478260869
$ws.Range("I23").Value = 156
$ws.Range("J23").Value = 149
$ws.Range("K23").Value = 4.697986577181
$ws.Range("L23").Value = 22.834645669291
$ws.Range("M23").Value = 52.941176470588

# Row 24
$ws.Range("C24").Value = 49
$ws.Range("D24").Value = 49
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 195
$ws.Range("G24").Value = 239
$ws.Range("H24").Value = -18.410041841004
$ws.Range("I24").Value = 1478
$ws.Range("J24").Value = 1713
$ws.Range("K24").Value = -13.718622300058
$ws.Range("L24").Value = 0.271370420624
$ws.Range("M24").Value = 55.907172995780

# Row 25
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 18
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 539
$ws.Range("J25").Value = 618
$ws.Range("K25").Value = -12.783171521035
$ws.Range("L25").Value = -2.531645569620
$ws.Range("M25").Value = -2

# Row 26
$ws.Range("D26").Value = '0'
$ws.Range("C26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = '***.*'
$ws.Range("C26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = -8.571428571428

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = 18.181818181818
$ws.Range("I27").Value = 70
$ws.Range("J27").Value = 53
$ws.Range("K27").Value = 32.075471698113
$ws.Range("L27").Value = -7.894736842105

# Row 28
$ws.Range("C28").Value = '0'
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -83.333333333333

# Row 29
$ws.Range("C29").Value = '0'
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -80

# Row 30
$ws.Range("D30").Value = 1
$ws.Range("F30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 160

